$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2746
$ws.Range("I28").Value = 994.4286
$ws.Range("J28").Value = 6833
$ws.Range("K28").Value = 994.4286
$ws.Range("L28").Value = 6833
$ws.Range("M28").Value = -509.4286
$ws.Range("N28").Value = -7803
$ws.Range("H37").Value = 1155
$ws.Range("I37").Value = 1155
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 3465
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -3339
$ws.Range("N37").ClearContents()
$ws.Range("H40").Value = 3459.2778
$ws.Range("I40").Value = 3730.1667
$ws.Range("J40").Value = 3323.8333
$ws.Range("K40").Value = 3730.1667
$ws.Range("L40").Value = 3323.8333
$ws.Range("M40").Value = -3555.1667
$ws.Range("N40").Value = -3673.8333
$ws.Range("H100").Value = 8665.666999999999
$ws.Range("I100").Value = 6000
$ws.Range("J100").Value = 9998.5
$ws.Range("K100").Value = 6000
$ws.Range("L100").Value = 9998.5
$ws.Range("M100").Value = -5459
$ws.Range("N100").Value = -11080.5
$ws.Range("H107").Value = 905.7143
$ws.Range("I107").Value = 974.63635
$ws.Range("J107").Value = 653
$ws.Range("K107").Value = 974.63635
$ws.Range("L107").Value = 653
$ws.Range("M107").Value = 945.36365
$ws.Range("N107").Value = -4493
$ws.Range("H113").Value = 3201.4285
$ws.Range("I113").Value = 1902
$ws.Range("J113").Value = 6450
$ws.Range("K113").Value = 1902
$ws.Range("L113").Value = 6450
$ws.Range("M113").Value = 1352
$ws.Range("N113").Value = -12958
$ws.Range("H137").Value = 31430.277
$ws.Range("I137").Value = 39224.81
$ws.Range("K137").Value = 117674.43
$ws.Range("M137").Value = -115124.43
$ws.Range("H138").Value = 340904.84
$ws.Range("I138").Value = 3522.577
$ws.Range("J138").Value = 523653.56
$ws.Range("K138").Value = 10567.731
$ws.Range("L138").Value = 1570960.68
$ws.Range("M138").Value = -5427.731
$ws.Range("N138").Value = -1581240.68

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2985.6274
$ws.Range("I32").Value = 2985.6274
$ws.Range("K32").Value = 2985.6274
$ws.Range("M32").Value = -2698.6274
$ws.Range("H61").Value = 5172.826
$ws.Range("I61").Value = 3560.1875
$ws.Range("K61").Value = 3560.1875
$ws.Range("M61").Value = -3348.1875
$ws.Range("H110").Value = 3170.1365
$ws.Range("I110").Value = 1884.6875
$ws.Range("K110").Value = 1884.6875
$ws.Range("M110").Value = 160.3125
$ws.Range("H122").Value = 4100.8237
$ws.Range("I122").Value = 3738.32
$ws.Range("J122").Value = 5107.778
$ws.Range("K122").Value = 11214.96
$ws.Range("L122").Value = 15323.334
$ws.Range("M122").Value = -8764.960000000001
$ws.Range("N122").Value = -20223.334
$ws.Range("H132").Value = 2670.525
$ws.Range("I132").Value = 1496.1936
$ws.Range("J132").Value = 6715.4443
$ws.Range("K132").Value = 4488.5808
$ws.Range("L132").Value = 20146.3329
$ws.Range("M132").Value = -1958.5808
$ws.Range("N132").Value = -25206.3329
$ws.Range("H136").Value = 5172.826
$ws.Range("I136").Value = 3560.1875
$ws.Range("K136").Value = 10680.5625
$ws.Range("M136").Value = -8130.5625

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 17860956
$ws.Range("I20").Value = 25004286
$ws.Range("J20").Value = 2625.875
$ws.Range("K20").Value = 25004286
$ws.Range("L20").Value = 2625.875
$ws.Range("M20").Value = -25004039
$ws.Range("N20").Value = -3119.875
$ws.Range("H86").Value = 3329.2222
$ws.Range("J86").Value = 3089
$ws.Range("L86").Value = 3089
$ws.Range("N86").Value = -5335
$ws.Range("H89").Value = 3329.2222
$ws.Range("J89").Value = 3089
$ws.Range("L89").Value = 15445
$ws.Range("N89").Value = -26677
$ws.Range("H99").Value = 79754.234
$ws.Range("I99").Value = 85775.414
$ws.Range("K99").Value = 85775.414
$ws.Range("M99").Value = -84277.414
$ws.Range("H105").Value = 14447156
$ws.Range("I105").Value = 835451
$ws.Range("J105").Value = 41670564
$ws.Range("K105").Value = 835451
$ws.Range("L105").Value = 41670564
$ws.Range("M105").Value = -833704
$ws.Range("N105").Value = -41674058
$ws.Range("H107").Value = 1273.875
$ws.Range("I107").Value = 1273.875
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1273.875
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 646.125
$ws.Range("N107").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H31").Value = 5283.2705
$ws.Range("I31").Value = 3814.4736
$ws.Range("J31").Value = 6833.6665
$ws.Range("K31").Value = 3814.4736
$ws.Range("L31").Value = 6833.6665
$ws.Range("M31").Value = -3519.4736
$ws.Range("N31").Value = -7423.6665
$ws.Range("H34").Value = 5283.2705
$ws.Range("I34").Value = 3814.4736
$ws.Range("J34").Value = 6833.6665
$ws.Range("K34").Value = 3814.4736
$ws.Range("L34").Value = 6833.6665
$ws.Range("M34").Value = -3612.4736
$ws.Range("N34").Value = -7237.6665
$ws.Range("H58").Value = 2745.4075
$ws.Range("I58").Value = 2197.0625
$ws.Range("K58").Value = 2197.0625
$ws.Range("M58").Value = -1994.0625
$ws.Range("H132").Value = 2211.5454
$ws.Range("I132").Value = 1981.1143
$ws.Range("J132").Value = 3107.6667
$ws.Range("K132").Value = 5943.3429
$ws.Range("L132").Value = 9323.000100000001
$ws.Range("M132").Value = -3413.3429
$ws.Range("N132").Value = -14383.0001
$ws.Range("H136").Value = 2745.4075
$ws.Range("I136").Value = 2197.0625
$ws.Range("K136").Value = 6591.1875
$ws.Range("M136").Value = -4041.1875
$ws.Range("H141").Value = 655949.7
$ws.Range("J141").Value = 722133.4399999999
$ws.Range("L141").Value = 722133.4399999999
$ws.Range("N141").Value = -732493.4399999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 372.8125
$ws.Range("I40").Value = 69
$ws.Range("J40").Value = 2499.5
$ws.Range("K40").Value = 276
$ws.Range("L40").Value = 9998
$ws.Range("M40").Value = -207
$ws.Range("N40").Value = -10136

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 17833
$ws.Range("I58").Value = 10000
$ws.Range("J58").Value = 21749.5
$ws.Range("K58").Value = 10000
$ws.Range("L58").Value = 21749.5
$ws.Range("M58").Value = -9723
$ws.Range("N58").Value = -22303.5
$ws.Range("H80").Value = 66669212
$ws.Range("I80").Value = 83335640
$ws.Range("J80").Value = 3500.3333
$ws.Range("K80").Value = 83335640
$ws.Range("L80").Value = 3500.3333
$ws.Range("M80").Value = -83334642
$ws.Range("N80").Value = -5496.3333
$ws.Range("H83").Value = 66669212
$ws.Range("I83").Value = 83335640
$ws.Range("J83").Value = 3500.3333
$ws.Range("K83").Value = 416678200
$ws.Range("L83").Value = 17501.6665
$ws.Range("M83").Value = -416673208
$ws.Range("N83").Value = -27485.6665
$ws.Range("H102").Value = 6705.5
$ws.Range("I102").Value = 741.25
$ws.Range("J102").Value = 18634
$ws.Range("K102").Value = 741.25
$ws.Range("L102").Value = 18634
$ws.Range("M102").Value = 880.75
$ws.Range("N102").Value = -21878
$ws.Range("H126").Value = 10497.833
$ws.Range("I126").Value = 2999.3333
$ws.Range("K126").Value = 8997.999899999999
$ws.Range("M126").Value = -6527.999899999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 100000000
$ws.Range("J2").Value = 100000000
$ws.Range("L2").Value = 100000000
$ws.Range("N2").Value = -100000224
$ws.Range("H61").Value = 4171
$ws.Range("I61").Value = 4004.3
$ws.Range("J61").Value = 5004.5
$ws.Range("K61").Value = 4004.3
$ws.Range("L61").Value = 5004.5
$ws.Range("M61").Value = -3802.3
$ws.Range("N61").Value = -5408.5
$ws.Range("H113").Value = 4171
$ws.Range("I113").Value = 4004.3
$ws.Range("J113").Value = 5004.5
$ws.Range("K113").Value = 4004.3
$ws.Range("L113").Value = 5004.5
$ws.Range("M113").Value = -1834.3
$ws.Range("N113").Value = -9344.5
$ws.Range("H122").Value = 3371.1428
$ws.Range("I122").Value = 3833.2222
$ws.Range("J122").Value = 2539.4
$ws.Range("K122").Value = 11499.6666
$ws.Range("L122").Value = 7618.200000000001
$ws.Range("M122").Value = -9049.6666
$ws.Range("N122").Value = -12518.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9261459
$ws.Range("I122").Value = 2264.2354
$ws.Range("J122").Value = 25002088
$ws.Range("K122").Value = 6792.706200000001
$ws.Range("L122").Value = 75006264
$ws.Range("M122").Value = -4342.706200000001
$ws.Range("N122").Value = -75011164
$ws.Range("H136").Value = 58828020
$ws.Range("I136").Value = 76924264
$ws.Range("K136").Value = 230772792
$ws.Range("M136").Value = -230770242

Write-Output "Applied all cell updates"